$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "BannedPaths" rule (row 35) is being replaced by a corrected
# "BannedPath" rule that now lives further down the table (just before the
# "AEM Rules:AEM-3" row) with an updated Severity of "Critical" and no Tags.

# 1) Remove the old "BannedPaths" row - this shifts rows 36-40 up to 35-39.
$ws.Rows(35).Delete()

# 2) Insert a new row in its new location (row 40, right before the row
#    that now holds "AEM Rules:AEM-3") and populate it.
$ws.Rows(40).Insert()
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"

# 3) Restore the active selection to match the saved workbook state.
$ws.Range("A37").Select()
